$wb = $excel.ActiveWorkbook

# Update the "zh-cn" sheet: Correspond Handoff/Handback Datetime values.
# Rows 2 and 4 originally shared the same timestamp text, so both rows
# need to be updated to keep them in sync (matching the shared-string
# table edit in the workbook).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 18:20:53"
$wsZhCn.Range("H2").Value = "2016-03-21 18:21:26"
$wsZhCn.Range("E4").Value = "2016-03-21 18:20:53"
$wsZhCn.Range("H4").Value = "2016-03-21 18:21:26"

# Update the "de-de" sheet: Correspond Handoff/Handback Datetime values.
# Same situation as above: rows 2 and 4 shared the same timestamp text.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 18:20:57"
$wsDeDe.Range("H2").Value = "2016-03-21 18:21:32"
$wsDeDe.Range("E4").Value = "2016-03-21 18:20:57"
$wsDeDe.Range("H4").Value = "2016-03-21 18:21:32"
